$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText.Replace(
    "✅ 1000 Bs = 4.74 = 18736.08 pesos`n✅ 18736.08 pesos = 4.7 = 931.16 Bs",
    "✅ 1000 Bs = 4.75 = 18779.1 pesos`n✅ 18779.1 pesos = 4.72 = 938.96 Bs"
)
$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 210.5
$wsTasas.Range("O10").Value = 3953
$wsTasas.Range("N12").Value = 3980
$wsTasas.Range("O12").Value = 199.001
